$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new match data (English Premier League, Man Utd vs Newcastle)
$ws.Range("A2").Value = "English Premier League"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "2025-12-26"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "17:00:00"
$ws.Range("D2").Value = "Man Utd"
$ws.Range("E2").Value = "Newcastle"
$ws.Range("F2").Value = 1.6
$ws.Range("G2").Value = 1.61
$ws.Range("H2").Value = 7.4
$ws.Range("I2").Value = 7.8
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 4.1
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 6.6
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 2.08
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 1.33
$ws.Range("S2").Value = 3.95
$ws.Range("T2").Value = 1.63
$ws.Range("U2").Value = 2.52
$ws.Range("V2").Value = 1.15
$ws.Range("W2").Value = 2.6
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 6.8
$ws.Range("AC2").Value = 5.6
$ws.Range("AD2").Value = 13.5
$ws.Range("AE2").Value = 46
$ws.Range("AF2").Value = 6.4
$ws.Range("AG2").Value = 7.4
$ws.Range("AH2").Value = 17
$ws.Range("AI2").Value = 55
$ws.Range("AJ2").Value = 17.5
$ws.Range("AK2").Value = 19
$ws.Range("AL2").Value = 44
$ws.Range("AM2").Value = 150
$ws.Range("AN2").Value = 25
$ws.Range("AO2").Value = 85

# Remove rows 3 and 4 (Saudi Professional League, and old English Premier League row)
$ws.Rows("3:4").Delete()
